$d = $word.ActiveDocument

# --- 1. Split the big paragraph: move "For each transaction..." into its own paragraph ---
$full = $d.Content.Text
$idx = $full.IndexOf("For each transaction between a ticket retailer")
$rng = $d.Range($idx, $idx)
$rng.InsertParagraphBefore()

# --- 2. Fix "Through its easy to use" -> " Through easy to use" (drop "its", add leading space) ---
$full = $d.Content.Text
$idx = $full.IndexOf("Through its easy to use interface")
$rng = $d.Range($idx, $idx + ("Through its easy to use interface").Length)
$rng.Text = " Through easy to use interface"

# --- 3. Fix "ticket offerors" -> "ticket retailers" ---
$full = $d.Content.Text
$idx = $full.IndexOf("offerors")
$rng = $d.Range($idx, $idx + ("offerors").Length)
$rng.Text = "retailers"

# --- 4. Fix "companies.." -> "companies, etc." ---
$full = $d.Content.Text
$idx = $full.IndexOf("companies..")
$rng = $d.Range($idx, $idx + ("companies..").Length)
$rng.Text = "companies, etc."

Write-Host "Final text:"
Write-Host $d.Content.Text
Write-Host "Paragraph count: $($d.Paragraphs.Count)"
